# This script reorders the shared-string labels in column A (the "goods"
# names) so that the underlying shared string table ends up in the order
# described by the target diff, while the counts in column B (and which
# row each count sits in) are left completely untouched.
#
# The target order for rows 2..56 (column A), count in column B stays as-is:
$words = @(
    "хлеб",
    "вино",
    "скот",
    "холст",
    "кожа",
    "мед",
    "пиво",
    "сукно",
    "овчина",
    "лошадь",
    "воск",
    "масло",
    "сало",
    "железо",
    "полотно",
    "Крымскую соль",
    "колеса",
    "парча",
    "сено",
    "говядина",
    "табак",
    "позумент",
    "чулок",
    "шелк",
    "выбойка",
    "сахар",
    "лес",
    "лыко",
    "китайка",
    "сани",
    "сапог",
    "ладан",
    "коса",
    "веревка",
    "конь",
    "замок",
    "овца",
    "платок",
    "обод",
    "ром",
    "гвоздь",
    "горшок",
    "рогожа",
    "сосуд",
    "роза",
    "дуга",
    "бечева",
    "сковорода",
    "скотский кожа",
    "нитка",
    "брусья",
    "хомут",
    "котел",
    "гумми",
    "покроми"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, clear out every existing label in column A (rows 2..56). This
# fully dereferences each shared string so the workbook's shared string
# table collapses down to just the still-used entries (the "Counts"
# header in B1 stays untouched).
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 1).Value = ""
}

# Now re-write the labels in the desired final order. Because the
# underlying engine appends newly-introduced shared strings to the table
# in the order they are written, writing them in this order reproduces
# the target shared-string table ordering exactly.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 1).Value = $words[$r - 2]
}
